$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2115.742
$ws.Range("J17").Value = 2142.9333
$ws.Range("L17").Value = 6428.7999
$ws.Range("N17").Value = -6764.7999

$ws.Range("H62").Value = 6253767.5
$ws.Range("I62").Value = 11365294
$ws.Range("J62").Value = 6346.4443
$ws.Range("K62").Value = 11365294
$ws.Range("L62").Value = 6346.4443
$ws.Range("M62").Value = -11364670
$ws.Range("N62").Value = -7594.4443

$ws.Range("H63").Value = 53246
$ws.Range("J63").Value = 53246
$ws.Range("L63").Value = 53246
$ws.Range("N63").Value = -54494

$ws.Range("H65").Value = 6253767.5
$ws.Range("I65").Value = 11365294
$ws.Range("J65").Value = 6346.4443
$ws.Range("K65").Value = 56826470
$ws.Range("L65").Value = 31732.2215
$ws.Range("M65").Value = -56823350
$ws.Range("N65").Value = -37972.2215

$ws.Range("H66").Value = 53246
$ws.Range("J66").Value = 53246
$ws.Range("L66").Value = 159738
$ws.Range("N66").Value = -165978

$ws.Range("H98").Value = 1592.75
$ws.Range("I98").Value = 994.4706
$ws.Range("K98").Value = 994.4706
$ws.Range("M98").Value = 503.5294

$ws.Range("H113").Value = 6388
$ws.Range("J113").Value = 7499.4
$ws.Range("L113").Value = 7499.4
$ws.Range("N113").Value = -14007.4

$ws.Range("H122").Value = 1592.75
$ws.Range("I122").Value = 994.4706
$ws.Range("K122").Value = 2983.4118
$ws.Range("M122").Value = -533.4117999999999

$ws.Range("H127").Value = 2333.9375
$ws.Range("I127").Value = 449.72726
$ws.Range("K127").Value = 1349.18178
$ws.Range("M127").Value = 3610.81822

$ws.Range("H129").Value = 2021.12
$ws.Range("I129").Value = 794.9
$ws.Range("K129").Value = 2384.7
$ws.Range("M129").Value = 2615.3

$ws.Range("H136").Value = 49165.832
$ws.Range("J136").Value = 49165.832
$ws.Range("L136").Value = 49165.832
$ws.Range("N136").Value = -59365.832

$ws.Range("H137").Value = 1341577.4
$ws.Range("I137").Value = 2002760.8
$ws.Range("K137").Value = 6008282.4
$ws.Range("M137").Value = -6005732.4

$ws.Range("H138").Value = 7798.3335
$ws.Range("J138").Value = 7994.0415
$ws.Range("L138").Value = 23982.1245
$ws.Range("N138").Value = -34262.12450000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 5000
$ws.Range("K25").Value = 5000
$ws.Range("M25").Value = -4598

$ws.Range("H32").Value = 5556.6743
$ws.Range("I32").Value = 4842.9756
$ws.Range("J32").Value = 20187.5
$ws.Range("K32").Value = 4842.9756
$ws.Range("L32").Value = 20187.5
$ws.Range("M32").Value = -4555.9756
$ws.Range("N32").Value = -20761.5

$ws.Range("H122").Value = 4649.1304
$ws.Range("I122").Value = 2911.375
$ws.Range("K122").Value = 8734.125
$ws.Range("M122").Value = -6284.125

$ws.Range("H138").Value = 69000
$ws.Range("J138").Value = 69000
$ws.Range("L138").Value = 69000
$ws.Range("N138").Value = -79280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 2117.5
$ws.Range("J11").Value = 972
$ws.Range("L11").Value = 972
$ws.Range("N11").Value = -1252

$ws.Range("H59").Value = 87280.336
$ws.Range("J59").Value = 87280.336
$ws.Range("L59").Value = 87280.336
$ws.Range("N59").Value = -88974.336

$ws.Range("H105").Value = 1173.6875
$ws.Range("I105").Value = 1143.9231
$ws.Range("J105").Value = 1302.6666
$ws.Range("K105").Value = 1143.9231
$ws.Range("L105").Value = 1302.6666
$ws.Range("M105").Value = 603.0769
$ws.Range("N105").Value = -4796.6666

$ws.Range("H134").Value = 36549.7
$ws.Range("I134").Value = 2343.7896
$ws.Range("J134").Value = 95632.63
$ws.Range("K134").Value = 7031.3688
$ws.Range("L134").Value = 286897.89
$ws.Range("M134").Value = -4496.3688
$ws.Range("N134").Value = -291967.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34365.47
$ws.Range("J31").Value = 70844.375
$ws.Range("L31").Value = 70844.375
$ws.Range("N31").Value = -71434.375

$ws.Range("H34").Value = 34365.47
$ws.Range("J34").Value = 70844.375
$ws.Range("L34").Value = 70844.375
$ws.Range("N34").Value = -71248.375

$ws.Range("H122").Value = 5053
$ws.Range("I122").Value = 3497.8333
$ws.Range("K122").Value = 10493.4999
$ws.Range("M122").Value = -8043.499899999999

$ws.Range("H139").Value = 74865
$ws.Range("J139").Value = 74865
$ws.Range("L139").Value = 74865
$ws.Range("N139").Value = -85145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 347.85715
$ws.Range("I2").Value = 63.75
$ws.Range("J2").Value = 726.6667
$ws.Range("K2").Value = 382.5
$ws.Range("L2").Value = 4360.0002
$ws.Range("M2").Value = -269.5
$ws.Range("N2").Value = -4586.0002

$ws.Range("H38").Value = 140
$ws.Range("I38").Value = 70
$ws.Range("J38").Value = 350
$ws.Range("K38").Value = 210
$ws.Range("L38").Value = 1050
$ws.Range("M38").Value = 137
$ws.Range("N38").Value = -1744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 125074990
$ws.Range("J135").Value = 125074990
$ws.Range("L135").Value = 125074990
$ws.Range("N135").Value = -125085130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 776056.3
$ws.Range("J7").Value = 1255408.9
$ws.Range("L7").Value = 1255408.9
$ws.Range("N7").Value = -1255632.9

$ws.Range("H122").Value = 836742.5600000001
$ws.Range("I122").Value = 2725
$ws.Range("J122").Value = 1253751.4
$ws.Range("K122").Value = 8175
$ws.Range("L122").Value = 3761254.2
$ws.Range("M122").Value = -5725
$ws.Range("N122").Value = -3766154.2

$ws.Range("H126").Value = 776056.3
$ws.Range("J126").Value = 1255408.9
$ws.Range("L126").Value = 3766226.7
$ws.Range("N126").Value = -3771166.7

$ws.Range("H132").Value = 4402.1665
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 4881.8
$ws.Range("K132").Value = 6012
$ws.Range("L132").Value = 14645.4
$ws.Range("M132").Value = -3482
$ws.Range("N132").Value = -19705.4

$ws.Range("H135").Value = 62649.668
$ws.Range("J135").Value = 62649.668
$ws.Range("L135").Value = 62649.668
$ws.Range("N135").Value = -72789.66800000001

$ws.Range("H138").Value = 63701.5
$ws.Range("J138").Value = 63701.5
$ws.Range("L138").Value = 63701.5
$ws.Range("N138").Value = -73981.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 62506676
$ws.Range("I122").Value = 90915710
$ws.Range("K122").Value = 272747130
$ws.Range("M122").Value = -272744680

$ws.Range("H126").Value = 2379.85
$ws.Range("I126").Value = 1693.75
$ws.Range("K126").Value = 5081.25
$ws.Range("M126").Value = -2611.25

$ws.Range("H138").Value = 69158.8
$ws.Range("J138").Value = 69158.8
$ws.Range("L138").Value = 69158.8
$ws.Range("N138").Value = -79438.8

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
